# Update the "Förändrad" (Changed) date in column C for every data row
# (rows 2-176) from 2023-11-13 (serial 45243) to 2023-11-14 (serial 45244).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
